# Auto-generated Excel COM-interop script applying the Durandal_Profits leve-profit refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 4566.4614
$ws.Range("I74").Value = 3260
$ws.Range("J74").Value = 4958.4
$ws.Range("K74").Value = 3260
$ws.Range("L74").Value = 4958.4
$ws.Range("M74").Value = -2324
$ws.Range("N74").Value = -6830.4

$ws.Range("H77").Value = 4566.4614
$ws.Range("I77").Value = 3260
$ws.Range("J77").Value = 4958.4
$ws.Range("K77").Value = 16300
$ws.Range("L77").Value = 24792
$ws.Range("M77").Value = -11620
$ws.Range("N77").Value = -34152

$ws.Range("H137").Value = 946.2222
$ws.Range("I137").Value = 843.3333
$ws.Range("J137").Value = 1357.7778
$ws.Range("K137").Value = 2529.9999
$ws.Range("L137").Value = 4073.3334
$ws.Range("M137").Value = 20.0001000000002
$ws.Range("N137").Value = -9173.3334

$ws.Range("H138").Value = 4054.1345
$ws.Range("I138").Value = 1691
$ws.Range("K138").Value = 5073
$ws.Range("M138").Value = 67

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2326.0303
$ws.Range("I61").Value = 2018.3334
$ws.Range("K61").Value = 2018.3334
$ws.Range("M61").Value = -1806.3334

$ws.Range("H63").Value = 4036.3635
$ws.Range("I63").Value = 2800
$ws.Range("J63").Value = 7333.3335
$ws.Range("K63").Value = 2800
$ws.Range("L63").Value = 7333.3335
$ws.Range("M63").Value = -2114
$ws.Range("N63").Value = -8705.333500000001

$ws.Range("H66").Value = 4036.3635
$ws.Range("I66").Value = 2800
$ws.Range("J66").Value = 7333.3335
$ws.Range("K66").Value = 14000
$ws.Range("L66").Value = 36666.6675
$ws.Range("M66").Value = -10568
$ws.Range("N66").Value = -43530.6675

$ws.Range("H74").Value = 778.0645
$ws.Range("I74").Value = 627.48
$ws.Range("J74").Value = 1405.5
$ws.Range("K74").Value = 627.48
$ws.Range("L74").Value = 1405.5
$ws.Range("M74").Value = 246.52
$ws.Range("N74").Value = -3153.5

$ws.Range("H77").Value = 778.0645
$ws.Range("I77").Value = 627.48
$ws.Range("J77").Value = 1405.5
$ws.Range("K77").Value = 3137.4
$ws.Range("L77").Value = 7027.5
$ws.Range("M77").Value = 1230.6
$ws.Range("N77").Value = -15763.5

$ws.Range("H136").Value = 2326.0303
$ws.Range("I136").Value = 2018.3334
$ws.Range("K136").Value = 6055.0002
$ws.Range("M136").Value = -3505.0002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1803793.2
$ws.Range("I31").Value = 2565973.2
$ws.Range("J31").Value = 2277
$ws.Range("K31").Value = 2565973.2
$ws.Range("L31").Value = 2277
$ws.Range("M31").Value = -2565678.2
$ws.Range("N31").Value = -2867

$ws.Range("H34").Value = 1803793.2
$ws.Range("I34").Value = 2565973.2
$ws.Range("J34").Value = 2277
$ws.Range("K34").Value = 2565973.2
$ws.Range("L34").Value = 2277
$ws.Range("M34").Value = -2565771.2
$ws.Range("N34").Value = -2681

$ws.Range("H58").Value = 1130.4814
$ws.Range("I58").Value = 1088.5625
$ws.Range("J58").Value = 1191.4546
$ws.Range("K58").Value = 1088.5625
$ws.Range("L58").Value = 1191.4546
$ws.Range("M58").Value = -885.5625
$ws.Range("N58").Value = -1597.4546

$ws.Range("H75").Value = 20500
$ws.Range("J75").Value = 20500
$ws.Range("L75").Value = 20500
$ws.Range("N75").Value = -22496

$ws.Range("H78").Value = 20500
$ws.Range("J78").Value = 20500
$ws.Range("L78").Value = 61500
$ws.Range("N78").Value = -71484

$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("N80").ClearContents()

$ws.Range("H81").Value = 16000
$ws.Range("J81").Value = 16000
$ws.Range("L81").Value = 16000
$ws.Range("N81").Value = -17996

$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("N83").ClearContents()

$ws.Range("H84").Value = 16000
$ws.Range("J84").Value = 16000
$ws.Range("L84").Value = 48000
$ws.Range("N84").Value = -57984

$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("N87").ClearContents()

$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("N90").ClearContents()

$ws.Range("H105").Value = 2231.3333
$ws.Range("I105").Value = 1462.3334
$ws.Range("J105").Value = 3000.3333
$ws.Range("K105").Value = 1462.3334
$ws.Range("L105").Value = 3000.3333
$ws.Range("M105").Value = 284.6666
$ws.Range("N105").Value = -6494.3333

$ws.Range("H136").Value = 1130.4814
$ws.Range("I136").Value = 1088.5625
$ws.Range("J136").Value = 1191.4546
$ws.Range("K136").Value = 3265.6875
$ws.Range("L136").Value = 3574.3638
$ws.Range("M136").Value = -715.6875
$ws.Range("N136").Value = -8674.363799999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1222.4286
$ws.Range("I68").Value = 845.4545000000001
$ws.Range("J68").Value = 1395.2084
$ws.Range("K68").Value = 2536.3635
$ws.Range("L68").Value = 4185.6252
$ws.Range("M68").Value = -1725.3635
$ws.Range("N68").Value = -5807.6252

$ws.Range("H71").Value = 1222.4286
$ws.Range("I71").Value = 845.4545000000001
$ws.Range("J71").Value = 1395.2084
$ws.Range("K71").Value = 7609.0905
$ws.Range("L71").Value = 12556.8756
$ws.Range("M71").Value = -3553.0905
$ws.Range("N71").Value = -20668.8756

$ws.Range("H121").Value = 507574.75
$ws.Range("I121").Value = 15150
$ws.Range("J121").Value = 999999.5
$ws.Range("K121").Value = 45450
$ws.Range("L121").Value = 2999998.5
$ws.Range("M121").Value = -44140
$ws.Range("N121").Value = -3002618.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 72535.69
$ws.Range("I80").Value = 187618.33
$ws.Range("J80").Value = 3486.1
$ws.Range("K80").Value = 187618.33
$ws.Range("L80").Value = 3486.1
$ws.Range("M80").Value = -186620.33
$ws.Range("N80").Value = -5482.1

$ws.Range("H83").Value = 72535.69
$ws.Range("I83").Value = 187618.33
$ws.Range("J83").Value = 3486.1
$ws.Range("K83").Value = 938091.6499999999
$ws.Range("L83").Value = 17430.5
$ws.Range("M83").Value = -933099.6499999999
$ws.Range("N83").Value = -27414.5

$ws.Range("H102").Value = 1352.88
$ws.Range("I102").Value = 1368.1052
$ws.Range("J102").Value = 1304.6666
$ws.Range("K102").Value = 1368.1052
$ws.Range("L102").Value = 1304.6666
$ws.Range("M102").Value = 253.8948
$ws.Range("N102").Value = -4548.6666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 7905.0835
$ws.Range("I81").Value = 3423.75
$ws.Range("J81").Value = 16867.75
$ws.Range("K81").Value = 6847.5
$ws.Range("L81").Value = 33735.5
$ws.Range("M81").Value = -5786.5
$ws.Range("N81").Value = -35857.5

$ws.Range("H84").Value = 7905.0835
$ws.Range("I84").Value = 3423.75
$ws.Range("J84").Value = 16867.75
$ws.Range("K84").Value = 34237.5
$ws.Range("L84").Value = 168677.5
$ws.Range("M84").Value = -28933.5
$ws.Range("N84").Value = -179285.5

$ws.Range("H113").Value = 462.44446
$ws.Range("I113").Value = 462.44446
$ws.Range("K113").Value = 1387.33338
$ws.Range("M113").Value = 782.66662
